$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hard-coded "order amount" values for rows 2-14 (as in the original exercise),
# computed as: if quantity < reorder threshold -> max amount - quantity,
# otherwise -> "Skip reorder"
$ws.Range("F2").Value = "Skip reorder"
$ws.Range("F3").Value = "Skip reorder"
$ws.Range("F4").Value = "Skip reorder"
$ws.Range("F5").Value = 1136
$ws.Range("F6").Value = 1880
$ws.Range("F7").Value = "Skip reorder"
$ws.Range("F8").Value = 176
$ws.Range("F9").Value = 188
$ws.Range("F10").Value = "Skip reorder"
$ws.Range("F11").Value = "Skip reorder"
$ws.Range("F12").Value = "Skip reorder"
$ws.Range("F13").Value = "Skip reorder"
$ws.Range("F14").Value = 751
